# Project Sample Project is saved.TEST Author: admin. Type: SAVE.
# Rules sheet, row 11 (the R40 rule row): column B changes from the
# rule-name string "R40" to the string "1".
#
# A leading apostrophe forces Excel to store the numeric-looking literal
# as text (matching the workbook's <c t="s"> shared-string cell) instead
# of silently auto-converting it to the number 1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B11").Value = "'1"
